# Deploy the implementation guide:
#  - bump the "Date" metadata value
#  - add a new "NORMAL" concept row to the Concepts sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date property value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-01-30T21:30:05+00:00"

# --- Concepts sheet: append a new concept row (Level, Code, Display, Definition) ---
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Force column A to text first (Level values in this table are stored as text,
# e.g. "1"), then populate the new row's values.
$wsConcepts.Range("A6").NumberFormat = "@"
$wsConcepts.Range("A6").Value = "1"
$wsConcepts.Range("B6").Value = "NORMAL"
$wsConcepts.Range("C6").Value = "Normal Specimen"
$wsConcepts.Range("D6").Value = "Normal specimen"

# Copy the formatting (border/fill/alignment) of the row above onto the new row
# so the new row matches the existing data rows' style.
$wsConcepts.Range("A5:D5").Copy()
$wsConcepts.Range("A6:D6").PasteSpecial(-4122)
